# Link sw readings to chapter r code
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K5: add ALDA chapter 7 R code link into the existing extensions.html#residual-structure entry
$ws.Range("K5").Value = "extensions.html#residual-structure; https://stats.idre.ucla.edu/r/examples/alda/r-applied-longitudinal-data-analysis-ch-7/; https://canvas.uoregon.edu/files/10667756/download?download_frd=1"

# K6: new cell linking to ALDA chapters 4 and 5 R code
$ws.Range("K6").Value = "https://stats.idre.ucla.edu/r/examples/alda/r-applied-longitudinal-data-analysis-ch-4/; https://stats.idre.ucla.edu/r/examples/alda/r-applied-longitudinal-data-analysis-ch-5/"

# K10: prepend ALDA chapter 6 R code link before the existing nonlinearity reading
$ws.Range("K10").Value = "https://stats.idre.ucla.edu/r/examples/alda/r-applied-longitudinal-data-analysis-ch-6/; https://www.sds.pub/nonlinearity.html"

# Row 6 grew taller to accommodate the new K6 content
$ws.Rows(6).RowHeight = 68

# Update view state to reflect where the author was working
$null = $ws.Range("K10").Select()
